$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-22"

# Update the header label for the "through" date column (column I, header row 1)
$ws.Range("I1").Value = "2022 (through 08-22)"

# Update September value (row 9) in the "through" date column
$ws.Range("I9").Value = 129

# Update Total row (row 14) in the "through" date column
$ws.Range("I14").Value = 1100
